$d = $word.ActiveDocument

# Locate the (currently empty) "No Spacing" paragraph that comes
# right after the "Moi bang la 1 class..." paragraph. That is the
# paragraph that receives the new sentence and the _GoBack bookmark.
$findRng = $d.Content
$foundIt = $findRng.Find.Execute("1 class", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
$anchorPara = $findRng.Paragraphs(1)
$targetPara = $anchorPara.Next()

$newText = "Flie câu hỏi là k đc xóa"

# Remove any existing _GoBack bookmark -- it will be recreated below,
# right after the text we are about to insert.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert the new sentence, followed by a one-character sentinel. We
# then place the (initially collapsed) bookmark just before that
# sentinel -- i.e. NOT at the very last character position of the
# paragraph's text, a position that this runtime's Bookmarks.Add
# cannot place reliably -- and finally delete the sentinel so the
# bookmark ends up collapsed right after the real text, directly
# before the paragraph mark.
$targetPara.Range.InsertBefore($newText + "Z")

$targetPara2 = $anchorPara.Next()
$paraEnd = $targetPara2.Range.End
$bmPos = $paraEnd - 2

$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($bmPos, $bmPos + 1)
$sentinelRange.Delete()
